# Odobrany slide s bezpecnym ukoncenim, zmena poradia
#
# 1) Delete the "Bezpecne ulozeni" slide (SlideID 263).
# 2) Move the picture-only slide (SlideID 259) so it now sits right
#    before the "Checklist pozadavku" slide (SlideID 266), i.e. just
#    after the "Logovani" slide (SlideID 265).
#
# Resulting slide order (by SlideID):
#   256, 257, 260, 261, 262, 265, 259, 266, 264

$p = $ppt.ActivePresentation

function Get-SlideById($pres, $id) {
    for ($i = 1; $i -le $pres.Slides.Count; $i++) {
        $slide = $pres.Slides.Item($i)
        if ($slide.SlideID -eq $id) {
            return $slide
        }
    }
    return $null
}

# 1) Remove the "Bezpecne ulozeni" slide entirely.
$toDelete = Get-SlideById $p 263
$toDelete.Delete()

# 2) Re-order: move the picture slide (id 259) to sit right after the
#    "Logovani" slide (id 265), i.e. directly before "Checklist" (id 266).
#    The picture slide currently sits *before* "Logovani", so moving it
#    to Logovani's current index lands it directly after Logovani once
#    the shift happens (MoveTo's argument is the slide's final index).
$logovani = Get-SlideById $p 265
$picture = Get-SlideById $p 259
$picture.MoveTo($logovani.SlideIndex)
